$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.684.35"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "2.475.34"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("E4").Value = "  -0.14%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "319.18"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.56%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "93.15"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +1.62%  "
$ws.Range("E7").Value = "  +1.25%  "
$ws.Range("E8").Value = "  -0.03%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "33.18"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +3.11%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0859"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +9.06%  "
$ws.Range("E12").Value = "  +0.58%  "
$ws.Range("D13").Value = "2.856.49"
$ws.Range("E13").Value = "  +0.28%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.91"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +1.14%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "15.79"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.49%  "
$ws.Range("D16").Value = "2.477.49"
$ws.Range("E16").Value = "  -1.10%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.789"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +3.00%  "
$ws.Range("D18").Value = "41.639.85"
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("D20").Value = "0.0₃0953"
$ws.Range("E20").Value = "  +0.60%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "71.28"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.14%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "11.33"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +2.55%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "239.71"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +1.54%  "
$ws.Range("E24").Value = "  +1.12%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "1.94"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +3.04%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  +0.63%  "
$ws.Range("E28").Value = "  +2.76%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "9.85"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +2.01%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "36.15"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +2.33%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "158.54"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +1.72%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "5.53"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +2.07%  "
$ws.Range("E33").Value = "  -0.12%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "2.59"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.69%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.0768"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +1.63%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "17.32"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +1.07%  "
$ws.Range("E37").Value = "  +5.62%  "
$ws.Range("E38").Value = "  +2.72%  "
$ws.Range("E39").Value = "  +2.28%  "
$ws.Range("E40").Value = "  +0.96%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "4.01"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("E42").Value = "  +5.03%  "
$ws.Range("D43").Value = "1.994.66"
$ws.Range("E43").Value = "  +2.23%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.0285"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +1.01%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "19.02"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +2.48%  "
$ws.Range("E46").Value = "  +2.74%  "
$ws.Range("E47").Value = "  +3.19%  "
$ws.Range("D48").Value = "2.713.00"
$ws.Range("E48").Value = "  +0.18%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "97.40"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +0.44%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "74.36"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +3.78%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "67.23"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.81%  "
